# Qatar Stars League workbook update (02-05-2024 20:28)
# The source feed re-ordered several same-kick-off-time fixtures; this
# shows up as pairs of data rows whose match/odds details (everything
# except the row id in column A, the Div in column C, and the Date in
# column D) were swapped between the two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each tuple is a pair of worksheet row numbers whose B:AB contents
# (id, HomeTeam, AwayTeam, scores, odds, P/L columns) need to be
# exchanged with one another.
$rowPairs = @(
    @(29, 30),
    @(60, 61),
    @(112, 113),
    @(114, 115),
    @(122, 123),
    @(125, 127)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:AB$r1")
    $range2 = $ws.Range("B$r2`:AB$r2")

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value = $values2
    $range2.Value = $values1
}
